$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving values that look like plain numbers must be forced to Text format
# first, otherwise Excel auto-converts them (e.g. "1.000" -> 1, "0.00001095" -> 1.095E-05)
# and the original decimal-point formatting of the price strings would be lost.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data values
$ws.Range("D2").Value = "28.271.13"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.807.73"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "338.24"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.4682"
$ws.Range("E7").Value = "  +22.01%  "
$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").Value = "  +11.99%  "
$ws.Range("D9").Value = "45.28"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "1.154"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").Value = "0.07642"
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("D12").Value = "22.48"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "6.332"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").Value = "7.451"
$ws.Range("E15").Value = "  +4.72%  "
$ws.Range("D16").Value = "1.807.84"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "0.00001095"
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").Value = "0.06722"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "81.97"
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "17.43"
$ws.Range("E21").Value = "  +4.11%  "
$ws.Range("D22").Value = "6.436"
$ws.Range("E22").Value = "  +4.03%  "
$ws.Range("D23").Value = "28.280.20"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").Value = "11.88"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "2.410"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "20.78"
$ws.Range("E26").Value = "  +4.66%  "
$ws.Range("D27").Value = "154.19"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "2.378"
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("D29").Value = "2.013.26"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("D30").Value = "133.04"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "1.259"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "4.031"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "0.09641"
$ws.Range("E33").Value = "  +9.31%  "
$ws.Range("D34").Value = "5.862"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "0.2262"
$ws.Range("E35").Value = "  +7.87%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06389"
$ws.Range("E36").Value = "  +4.24%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "12.13"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "0.02355"
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("D39").Value = "5.266"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").Value = "0.6649"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.239"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.496"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "8.267"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "14.22"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "0.6143"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "3.863"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "130.78"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("D49").Value = "2.039"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").Value = "1.177"
$ws.Range("E51").Value = "  +0.45%  "
